$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Conectores" text to "Conectores h" (A15)
$ws.Range("A15").Value = "Conectores h"

# Update unit values (column C) which drive the shared formulas in column D
$ws.Range("C5").Value = 8.5
$ws.Range("C7").Value = 24
$ws.Range("C8").Value = 7.8
$ws.Range("C14").Value = 35
$ws.Range("C15").Value = 3.5

# Update the selection shown in the sheet view
$ws.Range("A4").Select()
